$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to retain their text representation so that
# values such as "30.256.15" or "5.140" are not re-interpreted as numbers/dates
# and lose formatting (trailing zeros, multiple separators, etc.).
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.256.15"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "1.865.86"
$ws.Range("E3").Value = "  -2.02%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "235.09"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "0.4676"
$ws.Range("E7").Value = "  -1.26%  "
$ws.Range("D8").Value = "0.2839"
$ws.Range("E8").Value = "  -0.89%  "
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("D10").Value = "20.73"
$ws.Range("E10").Value = "  +4.75%  "
$ws.Range("D11").Value = "0.07879"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").Value = "96.88"
$ws.Range("E12").Value = "  -3.83%  "
$ws.Range("D13").Value = "1.871.94"
$ws.Range("E13").Value = "  -1.64%  "
$ws.Range("D14").Value = "5.140"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").Value = "0.6758"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "279.93"
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").Value = "30.274.05"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").Value = "0.9997"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "5.480"
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("D20").Value = "12.66"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "2.109.68"
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("D22").Value = "0.000007258"
$ws.Range("E22").Value = "  -3.19%  "
$ws.Range("D23").Value = "0.9995"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "6.179"
$ws.Range("E24").Value = "  -1.51%  "
$ws.Range("D25").Value = "9.304"
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("D26").Value = "164.77"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("D27").Value = "19.10"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").Value = "1.908"
$ws.Range("E28").Value = "  -6.11%  "
$ws.Range("D29").Value = "1.351"
$ws.Range("E29").Value = "  -2.09%  "
$ws.Range("D30").Value = "0.09635"
$ws.Range("E30").Value = "  -3.12%  "
$ws.Range("D31").Value = "4.411"
$ws.Range("E31").Value = "  -2.22%  "
$ws.Range("D32").Value = "1.471"
$ws.Range("E32").Value = "  -2.98%  "
$ws.Range("D33").Value = "4.105"
$ws.Range("E33").Value = "  -3.78%  "
$ws.Range("D34").Value = "0.04706"
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("D35").Value = "0.7038"
$ws.Range("E35").Value = "  -2.65%  "
$ws.Range("D36").Value = "1.101"
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("D37").Value = "2.711"
$ws.Range("E37").Value = "  -0.45%  "
$ws.Range("D38").Value = "0.01871"
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("D39").Value = "6.370"
$ws.Range("E39").Value = "  -5.58%  "
$ws.Range("D40").Value = "2.530"
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("D41").Value = "73.50"
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("D42").Value = "1.941"
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("D43").Value = "0.8474"
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("D44").Value = "0.4179"
$ws.Range("E44").Value = "  -2.17%  "
$ws.Range("D45").Value = "104.25"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").Value = "0.9991"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("D47").Value = "7.197"
$ws.Range("E47").Value = "  -2.70%  "
$ws.Range("D48").Value = "9.309"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").Value = "935.75"
$ws.Range("E49").Value = "  -5.57%  "
$ws.Range("D50").Value = "34.11"
$ws.Range("E50").Value = "  -1.38%  "
$ws.Range("D51").Value = "0.1136"
$ws.Range("E51").Value = "  -4.21%  "

# Restore the original (default) cell formatting now that the text values are set.
$dataRange.ClearFormats()
